$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(8, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(9, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(10, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(11, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(12, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(13, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(14, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(15, 10).Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Cells.Item(16, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(17, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(18, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(19, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(20, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(21, 10).Value = '[''Belgium'', ''Mexico'', ''Paraguay'']'
$ws.Cells.Item(22, 10).Value = '[''Portugal'', ''Poland'', ''Morocco'']'
$ws.Cells.Item(24, 10).Value = '[''Portugal'', ''Morocco'', ''England'']'
$ws.Cells.Item(25, 10).Value = '[''Poland'', ''Morocco'', ''England'']'
$ws.Cells.Item(26, 10).Value = '[''Poland'', ''Morocco'', ''England'']'
$ws.Cells.Item(27, 10).Value = '[''Poland'', ''Morocco'', ''England'']'
$ws.Cells.Item(28, 10).Value = '[''Poland'', ''Morocco'', ''England'']'
$ws.Cells.Item(29, 10).Value = '[''Poland'', ''Morocco'', ''England'']'
$ws.Cells.Item(243, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(244, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(245, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(246, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(247, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(248, 10).Value = '[''Argentina'', ''Netherlands'']'
$ws.Cells.Item(256, 10).Value = '[''Croatia'', ''Brazil'']'
$ws.Cells.Item(257, 10).Value = '[''Croatia'', ''Brazil'']'
$ws.Cells.Item(261, 10).Value = '[''Croatia'', ''Brazil'']'
$ws.Cells.Item(262, 10).Value = '[''Croatia'', ''Brazil'']'
$ws.Cells.Item(268, 10).Value = '[''South Korea'', ''Switzerland'']'
$ws.Cells.Item(269, 10).Value = '[''South Korea'', ''Switzerland'']'
$ws.Cells.Item(270, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(271, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(272, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(279, 10).Value = '[''Argentina'', ''Greece'']'
$ws.Cells.Item(285, 10).Value = '[''United States'', ''Slovenia'']'
$ws.Cells.Item(287, 10).Value = '[''England'', ''United States'']'
$ws.Cells.Item(289, 10).Value = '[''Ghana'', ''Germany'']'
$ws.Cells.Item(290, 10).Value = '[''Ghana'', ''Germany'']'
$ws.Cells.Item(291, 10).Value = '[''Ghana'', ''Germany'']'
$ws.Cells.Item(292, 10).Value = '[''Ghana'', ''Germany'']'
$ws.Cells.Item(307, 10).Value = '[''Portugal'', ''Brazil'']'
$ws.Cells.Item(308, 10).Value = '[''Portugal'', ''Brazil'']'
$ws.Cells.Item(309, 10).Value = '[''Portugal'', ''Brazil'']'
$ws.Cells.Item(310, 10).Value = '[''Portugal'', ''Brazil'']'
$ws.Cells.Item(331, 10).Value = '[''Italy'', ''Costa Rica'']'
$ws.Cells.Item(332, 10).Value = '[''Uruguay'', ''Costa Rica'']'
$ws.Cells.Item(335, 10).Value = '[''Greece'', ''Colombia'']'
$ws.Cells.Item(336, 10).Value = '[''Greece'', ''Colombia'']'
$ws.Cells.Item(337, 10).Value = '[''Greece'', ''Colombia'']'
$ws.Cells.Item(341, 10).Value = '[''Greece'', ''Colombia'']'
$ws.Cells.Item(353, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(354, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(355, 10).Value = '[''France'', ''Switzerland'']'
$ws.Cells.Item(361, 10).Value = '[''Belgium'', ''Algeria'']'
$ws.Cells.Item(362, 10).Value = '[''Belgium'', ''Russia'']'
$ws.Cells.Item(363, 10).Value = '[''Belgium'', ''Algeria'']'
$ws.Cells.Item(364, 10).Value = '[''Belgium'', ''Algeria'']'
$ws.Cells.Item(383, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(387, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(388, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(389, 10).Value = '[''Mexico'', ''Germany'']'
$ws.Cells.Item(390, 10).Value = '[''Sweden'', ''Mexico'']'
$ws.Cells.Item(391, 10).Value = '[''Sweden'', ''Mexico'']'
$ws.Cells.Item(392, 10).Value = '[''Sweden'', ''Mexico'']'
$ws.Cells.Item(393, 10).Value = '[''Sweden'', ''Mexico'']'
$ws.Cells.Item(394, 10).Value = '[''Sweden'', ''Mexico'']'
$ws.Cells.Item(396, 10).Value = '[''Brazil'', ''Switzerland'']'
$ws.Cells.Item(402, 10).Value = '[''Senegal'', ''Japan'']'
$ws.Cells.Item(405, 10).Value = '[''Belgium'', ''England'']'
$ws.Cells.Item(406, 10).Value = '[''Belgium'', ''England'']'
$ws.Cells.Item(407, 10).Value = '[''Belgium'', ''England'']'
$ws.Cells.Item(408, 10).Value = '[''Belgium'', ''England'']'
$ws.Cells.Item(409, 10).Value = '[''Belgium'', ''England'']'
$ws.Cells.Item(410, 10).Value = '[''Ecuador'', ''Netherlands'']'
$ws.Cells.Item(411, 10).Value = '[''Ecuador'', ''Netherlands'']'
$ws.Cells.Item(414, 10).Value = '[''Ecuador'', ''Netherlands'']'
$ws.Cells.Item(421, 10).Value = '[''France'', ''Australia'']'
$ws.Cells.Item(422, 10).Value = '[''Tunisia'', ''France'']'
$ws.Cells.Item(423, 10).Value = '[''France'', ''Australia'']'
$ws.Cells.Item(424, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(425, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(426, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(427, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(428, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(429, 10).Value = '[''Argentina'', ''Poland'']'
$ws.Cells.Item(430, 10).Value = '[''Croatia'', ''Morocco'']'
$ws.Cells.Item(431, 10).Value = '[''Croatia'', ''Morocco'']'
$ws.Cells.Item(432, 10).Value = '[''Croatia'', ''Morocco'']'
$ws.Cells.Item(433, 10).Value = '[''Croatia'', ''Morocco'']'
$ws.Cells.Item(446, 10).Value = '[''Uruguay'', ''Portugal'']'
$ws.Cells.Item(447, 10).Value = '[''Uruguay'', ''Portugal'']'
$ws.Cells.Item(448, 10).Value = '[''Uruguay'', ''Portugal'']'
$ws.Cells.Item(30, 10).Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Cells.Item(32, 10).Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Cells.Item(34, 10).Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Cells.Item(35, 10).Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Cells.Item(36, 10).Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Cells.Item(40, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(41, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(42, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(43, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(44, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(45, 10).Value = '[''Cameroon'', ''Argentina'', ''Soviet Union'']'
$ws.Cells.Item(46, 10).Value = '[''Cameroon'', ''Argentina'', ''Romania'']'
$ws.Cells.Item(47, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(48, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(49, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(50, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(51, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(52, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(53, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(54, 10).Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Cells.Item(55, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(56, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(57, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(58, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(59, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(60, 10).Value = '[''Czechoslovakia'', ''Italy'', ''Austria'']'
$ws.Cells.Item(61, 10).Value = '[''Scotland'', ''Brazil'', ''Costa Rica'']'
$ws.Cells.Item(62, 10).Value = '[''Scotland'', ''Sweden'', ''Brazil'']'
$ws.Cells.Item(63, 10).Value = '[''Scotland'', ''Brazil'', ''Costa Rica'']'
$ws.Cells.Item(64, 10).Value = '[''Scotland'', ''Brazil'', ''Costa Rica'']'
$ws.Cells.Item(65, 10).Value = '[''Scotland'', ''Brazil'', ''Costa Rica'']'
$ws.Cells.Item(66, 10).Value = '[''Belgium'', ''Uruguay'', ''Spain'']'
$ws.Cells.Item(68, 10).Value = '[''Belgium'', ''Uruguay'', ''Spain'']'
$ws.Cells.Item(72, 10).Value = '[''Egypt'', ''England'', ''Netherlands'']'
$ws.Cells.Item(73, 10).Value = '[''Egypt'', ''England'', ''Netherlands'']'
$ws.Cells.Item(75, 10).Value = '[''Romania'', ''Switzerland'', ''United States'']'
$ws.Cells.Item(76, 10).Value = '[''Romania'', ''Switzerland'', ''United States'']'
$ws.Cells.Item(77, 10).Value = '[''Romania'', ''Switzerland'', ''United States'']'
$ws.Cells.Item(78, 10).Value = '[''Romania'', ''Switzerland'', ''United States'']'
$ws.Cells.Item(79, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(80, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(81, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(82, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(83, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(84, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(85, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(86, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(87, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(88, 10).Value = '[''South Korea'', ''Spain'', ''Germany'']'
$ws.Cells.Item(89, 10).Value = '[''Ireland'', ''Italy'', ''Mexico'']'
$ws.Cells.Item(90, 10).Value = '[''Ireland'', ''Italy'', ''Norway'']'
$ws.Cells.Item(91, 10).Value = '[''Ireland'', ''Italy'', ''Mexico'']'
$ws.Cells.Item(92, 10).Value = '[''Cameroon'', ''Sweden'', ''Brazil'']'
$ws.Cells.Item(93, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(94, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(95, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(96, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(97, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(98, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(99, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(100, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(101, 10).Value = '[''Sweden'', ''Brazil'', ''Russia'']'
$ws.Cells.Item(102, 10).Value = '[''Belgium'', ''Saudi Arabia'', ''Netherlands'']'
$ws.Cells.Item(103, 10).Value = '[''Belgium'', ''Saudi Arabia'', ''Netherlands'']'
$ws.Cells.Item(104, 10).Value = '[''Belgium'', ''Saudi Arabia'', ''Netherlands'']'
$ws.Cells.Item(105, 10).Value = '[''Belgium'', ''Saudi Arabia'', ''Netherlands'']'
$ws.Cells.Item(106, 10).Value = '[''Belgium'', ''Saudi Arabia'', ''Netherlands'']'
$ws.Cells.Item(109, 10).Value = '[''Argentina'', ''Bulgaria'', ''Nigeria'']'
$ws.Cells.Item(110, 10).Value = '[''Argentina'', ''Bulgaria'', ''Nigeria'']'
$ws.Cells.Item(111, 10).Value = '[''Argentina'', ''Bulgaria'', ''Nigeria'']'
$ws.Cells.Item(112, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(113, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(114, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(115, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(116, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(117, 10).Value = '[''Italy'', ''Chile'']'
$ws.Cells.Item(145, 10).Value = '[''Mexico'', ''Netherlands'']'
$ws.Cells.Item(146, 10).Value = '[''Mexico'', ''Netherlands'']'
$ws.Cells.Item(147, 10).Value = '[''Belgium'', ''Netherlands'']'
$ws.Cells.Item(148, 10).Value = '[''Belgium'', ''Netherlands'']'
$ws.Cells.Item(149, 10).Value = '[''Mexico'', ''Netherlands'']'
$ws.Cells.Item(150, 10).Value = '[''Mexico'', ''Netherlands'']'
$ws.Cells.Item(151, 10).Value = '[''Mexico'', ''Netherlands'']'
$ws.Cells.Item(152, 10).Value = '[''Yugoslavia'', ''Germany'']'
$ws.Cells.Item(153, 10).Value = '[''Yugoslavia'', ''Germany'']'
$ws.Cells.Item(154, 10).Value = '[''Yugoslavia'', ''Germany'']'
$ws.Cells.Item(155, 10).Value = '[''Yugoslavia'', ''Germany'']'
$ws.Cells.Item(156, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(157, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(158, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(159, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(160, 10).Value = '[''Argentina'', ''Croatia'']'
$ws.Cells.Item(166, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(167, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(168, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(169, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(170, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(171, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(172, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(173, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(174, 10).Value = '[''Denmark'', ''Senegal'']'
$ws.Cells.Item(175, 10).Value = '[''Cameroon'', ''Germany'']'
$ws.Cells.Item(176, 10).Value = '[''Ireland'', ''Germany'']'
$ws.Cells.Item(177, 10).Value = '[''Ireland'', ''Germany'']'
$ws.Cells.Item(178, 10).Value = '[''Ireland'', ''Germany'']'
$ws.Cells.Item(179, 10).Value = '[''Ireland'', ''Germany'']'
$ws.Cells.Item(180, 10).Value = '[''Ireland'', ''Germany'']'
$ws.Cells.Item(181, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(182, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(183, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(197, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(198, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(199, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(200, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(202, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(203, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(204, 10).Value = '[''Turkey'', ''Brazil'']'
$ws.Cells.Item(209, 10).Value = '[''Russia'', ''Japan'']'
$ws.Cells.Item(210, 10).Value = '[''Belgium'', ''Japan'']'
$ws.Cells.Item(211, 10).Value = '[''Belgium'', ''Japan'']'
$ws.Cells.Item(212, 10).Value = '[''Russia'', ''Japan'']'
$ws.Cells.Item(213, 10).Value = '[''Russia'', ''Japan'']'
$ws.Cells.Item(214, 10).Value = '[''Belgium'', ''Japan'']'
$ws.Cells.Item(215, 10).Value = '[''Belgium'', ''Japan'']'
$ws.Cells.Item(216, 10).Value = '[''Belgium'', ''Japan'']'
$ws.Cells.Item(217, 10).Value = '[''South Korea'', ''United States'']'
$ws.Cells.Item(221, 10).Value = '[''South Korea'', ''United States'']'
$ws.Cells.Item(222, 10).Value = '[''South Korea'', ''United States'']'
$ws.Cells.Item(223, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(224, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(225, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(226, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(227, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(228, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(229, 10).Value = '[''Ecuador'', ''Germany'']'
$ws.Cells.Item(230, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(231, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(232, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(233, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(234, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(235, 10).Value = '[''Sweden'', ''England'']'
$ws.Cells.Item(236, 10).Value = '[''Sweden'', ''England'']'
